$d = $word.ActiveDocument

# Helper: replace a paragraph's visible text while keeping its leading empty
# run and its <w:pPr> (paragraph properties) intact, by round-tripping
# through InsertXML on the paragraph's own Range (collapsing adjacent runs
# with identical/empty formatting is avoided this way).
function Set-BulletParaText($index, $newText) {
    $p = $d.Paragraphs.Item($index)
    $r = $p.Range
    $pkg = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?>' +
           '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
           '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
           '<w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/><w:spacing w:line="240" w:lineRule="auto"/><w:ind w:left="720"/></w:pPr>' +
           '<w:r/><w:r><w:t>' + $newText + '</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
    [void]$r.InsertXML($pkg)
}

# 1. Title heading (appears both as the H1 heading and again, bolded, near
#    the end of the document) - plain text substitution, no run-merging
#    concerns here since the two occurrences keep their own run shapes.
$d.Content.Find.Execute(
    "Play Brazil Bomba for Free - Yggdrasil's Carnival-Themed Slot", $true, $false, $false, $false, $false,
    $true, 1, $false, "Play Brazil Bomba Slot for Free", 2) | Out-Null

# 2. "What we like" bullet list (paragraphs 44-47)
Set-BulletParaText 44 "Themed after the famous carnival in Rio de Janeiro"
Set-BulletParaText 45 "Drop mechanism adds excitement to gameplay"
Set-BulletParaText 46 "Multiplier system for big wins"
Set-BulletParaText 47 "Excellent graphics and audio create an immersive experience"

# 3. "What we don't like" bullet list (paragraphs 49-50)
Set-BulletParaText 49 "No progressive jackpot"
Set-BulletParaText 50 "Limited number of special features"

# 4. Closing italic summary paragraph
$d.Content.Find.Execute(
    "Experience the excitement of Rio's carnival with Brazil Bomba, Yggdrasil's online slot game. Play for free and activate the multiplier system for big wins!",
    $true, $false, $false, $false, $false,
    $true, 1, $false, "Read our review of Brazil Bomba, a themed online slot with excellent graphics and exciting gameplay. Play for free and win big!", 2) | Out-Null
